# Update TPM-derived LR-pair statistics (Anxa2-Tlr2) in columns G:T
# for rows 2-10, per the "update scripts wuth new tpm" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 32.09557633333333
$ws.Range("H2").Value = 96.28672900000001
$ws.Range("I2").Value = 0.1656600924295661
$ws.Range("J2").Value = 0.1656600924295661
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1465046666666667
$ws.Range("N2").Value = 0.439514
$ws.Range("O2").Value = 0.07745172725947863
$ws.Range("P2").Value = 0.07745172725947864
$ws.Range("Q2").Value = 4.702151712189556
$ws.Range("R2").Value = 42.319365409706
$ws.Range("S2").Value = 0.01283066029663477
$ws.Range("T2").Value = 0.01283066029663477
$ws.Range("G3").Value = 32.09557633333333
$ws.Range("H3").Value = 96.28672900000001
$ws.Range("I3").Value = 0.1656600924295661
$ws.Range("J3").Value = 0.1656600924295661
$ws.Range("N3").Value = 4.707498
$ws.Range("O3").Value = 0.8295614045753745
$ws.Range("P3").Value = 0.8295614045753745
$ws.Range("Q3").Value = 50.36328713267134
$ws.Range("R3").Value = 453.2695841940421
$ws.Range("S3").Value = 0.1374252189579572
$ws.Range("T3").Value = 0.1374252189579572
$ws.Range("G4").Value = 32.09557633333333
$ws.Range("H4").Value = 96.28672900000001
$ws.Range("I4").Value = 0.1656600924295661
$ws.Range("J4").Value = 0.1656600924295661
$ws.Range("O4").Value = 0.09298686816514685
$ws.Range("P4").Value = 0.09298686816514684
$ws.Range("Q4").Value = 5.645301619795445
$ws.Range("R4").Value = 50.80771457815901
$ws.Range("S4").Value = 0.0154042131749741
$ws.Range("T4").Value = 0.0154042131749741
$ws.Range("I5").Value = 0.6938590312037638
$ws.Range("J5").Value = 0.6938590312037638
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1465046666666667
$ws.Range("N5").Value = 0.439514
$ws.Range("O5").Value = 0.07745172725947863
$ws.Range("P5").Value = 0.07745172725947864
$ws.Range("Q5").Value = 19.69472782336
$ws.Range("R5").Value = 177.25255041024
$ws.Range("S5").Value = 0.05374058044131998
$ws.Range("T5").Value = 0.05374058044132
$ws.Range("I6").Value = 0.6938590312037638
$ws.Range("J6").Value = 0.6938590312037638
$ws.Range("N6").Value = 4.707498
$ws.Range("O6").Value = 0.8295614045753745
$ws.Range("P6").Value = 0.8295614045753745
$ws.Range("S6").Value = 0.5755986725027029
$ws.Range("T6").Value = 0.5755986725027029
$ws.Range("I7").Value = 0.6938590312037638
$ws.Range("J7").Value = 0.6938590312037638
$ws.Range("O7").Value = 0.09298686816514685
$ws.Range("P7").Value = 0.09298686816514684
$ws.Range("S7").Value = 0.0645197782597409
$ws.Range("T7").Value = 0.0645197782597409
$ws.Range("I8").Value = 0.1404808763666701
$ws.Range("J8").Value = 0.1404808763666701
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1465046666666667
$ws.Range("N8").Value = 0.439514
$ws.Range("O8").Value = 0.07745172725947863
$ws.Range("P8").Value = 0.07745172725947864
$ws.Range("Q8").Value = 3.987456385238223
$ws.Range("R8").Value = 35.88710746714401
$ws.Range("S8").Value = 0.01088048652152387
$ws.Range("T8").Value = 0.01088048652152387
$ws.Range("I9").Value = 0.1404808763666701
$ws.Range("J9").Value = 0.1404808763666701
$ws.Range("N9").Value = 4.707498
$ws.Range("O9").Value = 0.8295614045753745
$ws.Range("P9").Value = 0.8295614045753745
$ws.Range("Q9").Value = 42.70840737404534
$ws.Range("S9").Value = 0.1165375131147144
$ws.Range("T9").Value = 0.1165375131147144
$ws.Range("I10").Value = 0.1404808763666701
$ws.Range("J10").Value = 0.1404808763666701
$ws.Range("O10").Value = 0.09298686816514685
$ws.Range("P10").Value = 0.09298686816514684
$ws.Range("S10").Value = 0.01306287673043185
$ws.Range("T10").Value = 0.01306287673043185
